$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 5: consolidate C5's value into B5, drop old B5 (=1) and C5 ---
$ws.Cells.Item(5, 2).Value = 0.058
$ws.Cells.Item(5, 3).ClearContents()

# --- Row 9: shift values left by one column, drop E9 ---
$ws.Cells.Item(9, 2).Value = 0.02
$ws.Cells.Item(9, 3).Value = 0.045
$ws.Cells.Item(9, 4).Value = 0.09
$ws.Cells.Item(9, 5).ClearContents()

# --- Insert 5 new parameter rows before the existing row 57 ---
$ws.Range("A57:A61").EntireRow.Insert()

$ws.Cells.Item(57, 1).Value = "odds_homebirth"
$ws.Cells.Item(57, 2).Value = 0.5

$ws.Cells.Item(58, 1).Value = "or_homebirth_unmarried"
$ws.Cells.Item(58, 2).Value = 1.83

$ws.Cells.Item(59, 1).Value = "or_homebirth_wealth_4"
$ws.Cells.Item(59, 2).Value = 0.51

$ws.Cells.Item(60, 1).Value = "or_homebirth_wealth_5"
$ws.Cells.Item(60, 2).Value = 0.43

$ws.Cells.Item(61, 1).Value = "or_homebirth_urban"
$ws.Cells.Item(61, 2).Value = 0.39
# Row 61's A cell keeps the default (no explicit) style, unlike A57:A60
$ws.Cells.Item(61, 1).Style = "Normal"

# --- Update the view: scroll down and move the selection ---
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("F59").Select()
